$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1821
$ws1.Range("F8").Value = 67
$ws1.Range("F18").Value = 5201
$ws1.Range("F19").Value = 63
$ws1.Range("F20").Value = 849
$ws1.Range("F22").Value = 2306
$ws1.Range("F24").Value = 36
$ws1.Range("F25").Value = 2150

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1821
$ws4.Range("F8").Value = 67
$ws4.Range("F18").Value = 5201
$ws4.Range("F20").Value = 63
$ws4.Range("F22").Value = 849
$ws4.Range("F24").Value = 2306
$ws4.Range("F27").Value = 36
$ws4.Range("F28").Value = 2150
